$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing quota values (column C) ---
$ws.Cells.Item(4, 3).Value = 40   # Deep Learning for Search: 30 -> 40
$ws.Cells.Item(5, 3).Value = 40   # MLOps Engineering: 30 -> 40
$ws.Cells.Item(7, 3).Value = 40   # Introduction to Neuroscience: 30 -> 40
$ws.Cells.Item(8, 3).Value = 40   # Introduction to Prototyping: 30 -> 40
$ws.Cells.Item(9, 3).Value = 40   # Architecture of Computing Devices: 30 -> 40
$ws.Cells.Item(10, 3).Value = 28  # Real-Time Scheduling Theory: 30 -> 28
$ws.Cells.Item(11, 3).Value = 40  # Introduction to Robotics Operating System...: 30 -> 40
$ws.Cells.Item(12, 3).Value = 40  # Lambda-Calculus, Algebra, ...: 35 -> 40

# --- Append new course rows ---
$ws.Cells.Item(13, 1).Value = 12
$ws.Cells.Item(13, 2).Value = "Cross-platform Mobile Development with Flutter"
$ws.Cells.Item(13, 2).Font.Size = 10
$ws.Cells.Item(13, 3).Value = 40

$ws.Cells.Item(14, 1).Value = 13
$ws.Cells.Item(14, 2).Value = "Advanced Programming in C/C++ (Russian only)"
$ws.Cells.Item(14, 2).Font.Size = 10
$ws.Cells.Item(14, 3).Value = 40

$ws.Cells.Item(15, 1).Value = 14
$ws.Cells.Item(15, 2).Value = "Introduction to Mechanical Engineering"
$ws.Cells.Item(15, 2).Font.Size = 10
$ws.Cells.Item(15, 3).Value = 40

$ws.Cells.Item(16, 1).Value = 15
$ws.Cells.Item(16, 2).Value = "Introduction to Electronic and Logic Circuits"
$ws.Cells.Item(16, 2).Font.Size = 10
$ws.Cells.Item(16, 2).Interior.Color = 13421812
$ws.Cells.Item(16, 3).Value = 40

# --- Update selection to reflect the new active cell ---
$ws.Range("P17").Select() | Out-Null
